# Weekly fruit/hortaliza price-sheet refresh.
# The source rows for "Mapocho Venta Directa de Santiago - Ají" were
# reshuffled (row 6 is untouched); update every data row (2-15) in place
# to hold its new reported values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44460
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 95000
$ws.Range("L2").Value = 95000
$ws.Range("M2").Value = 95000
$ws.Range("P2").Value = 3800

# Row 3
$ws.Range("D3").Value = 44193
$ws.Range("H3").Value = 'Americana (o)'
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 46000
$ws.Range("L3").Value = 46000
$ws.Range("M3").Value = 46000
$ws.Range("P3").Value = 3067

# Row 4
$ws.Range("D4").Value = 44221
$ws.Range("J4").Value = 22
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24545
$ws.Range("P4").Value = 982

# Row 5
$ws.Range("D5").Value = 44421

# Row 7
$ws.Range("D7").Value = 44319
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 30000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 30000
$ws.Range("P7").Value = 1200

# Row 8
$ws.Range("D8").Value = 44326
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = 30000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 30000
$ws.Range("P8").Value = 1200

# Row 9
$ws.Range("D9").Value = 44449
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 80000
$ws.Range("L9").Value = 80000
$ws.Range("M9").Value = 80000
$ws.Range("N9").Value = '$/caja 25 kilos'
$ws.Range("P9").Value = 3200
$ws.Range("Q9").Value = 25

# Row 10
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 75000
$ws.Range("L10").Value = 75000
$ws.Range("M10").Value = 75000
$ws.Range("N10").Value = '$/caja 15 kilos'
$ws.Range("P10").Value = 5000
$ws.Range("Q10").Value = 15

# Row 11
$ws.Range("D11").Value = 44425
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 15
$ws.Range("N11").Value = '$/caja 25 kilos'
$ws.Range("P11").Value = 3000
$ws.Range("Q11").Value = 25

# Row 12
$ws.Range("D12").Value = 44474
$ws.Range("J12").Value = 18
$ws.Range("K12").Value = 100000
$ws.Range("L12").Value = 100000
$ws.Range("M12").Value = 100000
$ws.Range("P12").Value = 4000

# Row 13
$ws.Range("D13").Value = 44446
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 78000
$ws.Range("L13").Value = 78000
$ws.Range("M13").Value = 78000
$ws.Range("P13").Value = 3120

# Row 14
$ws.Range("D14").Value = 44446
$ws.Range("H14").Value = 'Inferno'
$ws.Range("J14").Value = 4
$ws.Range("K14").Value = 80000
$ws.Range("L14").Value = 80000
$ws.Range("M14").Value = 80000
$ws.Range("N14").Value = '$/caja 15 kilos'
$ws.Range("P14").Value = 5333
$ws.Range("Q14").Value = 15

# Row 15
$ws.Range("D15").Value = 44340
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 35000
$ws.Range("L15").Value = 35000
$ws.Range("M15").Value = 35000
$ws.Range("P15").Value = 1400
